$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# The document ends with a single paragraph that holds both the visible
# text "Student Counseling Center" AND the (normally hidden) "_GoBack"
# bookmark. The new content must be inserted between that text and the
# bookmark, i.e. the bookmark has to remain the very last thing in the
# document, now living alone in its own trailing paragraph.
#
# Plain-text InsertBefore() (using a carriage return for paragraph
# breaks) reliably keeps a trailing bookmark anchored to the end of the
# inserted text, so the whole block is built that way first (using
# placeholder marker text), and then patched up:
#   * blank-placeholder paragraphs get their marker character deleted
#     with Range.Delete() (which - unlike Range.Text = "" - collapses
#     the paragraph down to a clean run-less <w:p/>);
#   * the two heading paragraphs ("Danley Hall:" / "Linscheid Library:")
#     get their marker text cleared and replaced with the real OOXML
#     (bold/underline/page-break run, spell-check-exception runs).
# ----------------------------------------------------------------------

$cr = [string][char]13
$enDash = [char]0x2013

$items = @(
    "<<BLANK1>>",
    "<<DANLEY_HALL>>",
    "101- IT HelpDesk",
    "125- IT Director",
    "120- Atrium",
    "201 $enDash Presidents Suite",
    "300- Marketing",
    "<<BLANK2>>",
    "<<LINSCHEID_LIBRARY>>",
    "Digital Humanities Lab",
    "201 Media Services",
    "401 Library Directors Office"
)

$block = $cr + ($items -join $cr) + $cr

# Insertion point: immediately at the current document end (i.e. right
# before the trailing _GoBack bookmark).
$tail = $d.Range($d.Content.End, $d.Content.End)
$tail.InsertBefore($block)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Clear-Marker([string]$marker) {
    $content = $d.Content
    $content.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $target = $d.Range($content.Start, $content.End)
    $target.Delete()
    return $target
}

# --- Remove the blank-paragraph placeholders ---------------------------
Clear-Marker("<<BLANK1>>") | Out-Null
Clear-Marker("<<BLANK2>>") | Out-Null

# --- Patch the "Danley Hall:" heading paragraph -------------------------
$target = Clear-Marker("<<DANLEY_HALL>>")
$danleyXml = "<w:p $wNs><w:r><w:rPr><w:b/><w:u w:val=`"single`"/></w:rPr><w:lastRenderedPageBreak/><w:t>Danley Hall</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>"
$target.InsertXML($danleyXml)

# --- Patch the "Linscheid Library:" heading paragraph -------------------
$target = Clear-Marker("<<LINSCHEID_LIBRARY>>")
$linscheidXml = "<w:p $wNs><w:proofErr w:type=`"spellStart`"/><w:r><w:t>Linscheid</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> Library:</w:t></w:r></w:p>"
$target.InsertXML($linscheidXml)
